# feat: tambah kolom kab_kota di template xlsx
#
# Inserts a new column before column A on Sheet1 and labels its header
# "kab_kota" (a new shared string). All existing columns (jenjang, npsn,
# nama_sekolah, orang_tua, nama_peserta, tempat_lahir, tanggal_lahir,
# tahun_lulus, nomor_ujian, nomor_ijazah) shift one column to the right,
# and the cell selection moves to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column at position A; existing columns A:J become B:K.
$ws.Columns("A:A").Insert()

# Set the header text for the freshly inserted column.
$ws.Range("A1").Value = "kab_kota"

# Move the active selection to A2, matching the saved selection state.
$ws.Range("A2").Select()
